$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to text format so numeric-looking strings
# (e.g. "314.53", "0.9973") are stored as text, matching the source data which
# uses free-form price strings (including multi-dot forms like "24.730.53").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.730.53'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '1.705.67'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("D4").Value = '0.9973'
$ws.Range("E4").Value = '  -0.83%  '
$ws.Range("D5").Value = '314.53'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").Value = '0.9982'
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("D7").Value = '0.3978'
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").Value = '0.4064'
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("D9").Value = '0.9973'
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").Value = '1.516'
$ws.Range("E10").Value = '  +7.15%  '
$ws.Range("D11").Value = '53.14'
$ws.Range("E11").Value = '  +9.49%  '
$ws.Range("D12").Value = '0.08822'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '7.337'
$ws.Range("E13").Value = '  +11.02%  '
$ws.Range("D14").Value = '23.49'
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("D15").Value = '0.00001328'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '7.571'
$ws.Range("E16").Value = '  +4.70%  '
$ws.Range("D17").Value = '1.699.74'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").Value = '100.98'
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").Value = '0.07129'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("D20").Value = '19.60'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").Value = '6.780'
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").Value = '0.9978'
$ws.Range("D23").Value = '14.22'
$ws.Range("E23").Value = '  +2.08%  '
$ws.Range("D24").Value = '24.710.77'
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("D25").Value = '3.021'
$ws.Range("E25").Value = '  +7.98%  '
$ws.Range("D26").Value = '2.310'
$ws.Range("D27").Value = '22.48'
$ws.Range("E27").Value = '  +1.37%  '
$ws.Range("D28").Value = '159.80'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = '5.127'
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("D30").Value = '133.62'
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("D31").Value = '7.331'
$ws.Range("E31").Value = '  +25.04%  '
$ws.Range("D32").Value = '1.886.99'
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("D33").Value = '1.091'
$ws.Range("E33").Value = '  -7.11%  '
$ws.Range("D34").Value = '0.08714'
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("D35").Value = '7.296'
$ws.Range("E35").Value = '  +17.99%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '11.07'
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("B37").Value = 'WEMIXTOKEN'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '1.955'
$ws.Range("E37").Value = '  +3.89%  '
$ws.Range("D38").Value = '0.2735'
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").Value = '14.92'
$ws.Range("E39").Value = '  -2.82%  '
$ws.Range("D40").Value = '0.02792'
$ws.Range("E40").Value = '  +9.51%  '
$ws.Range("D41").Value = '0.09019'
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("D42").Value = '1.484'
$ws.Range("E42").Value = '  +1.95%  '
$ws.Range("D43").Value = '0.7711'
$ws.Range("E43").Value = '  +1.91%  '
$ws.Range("D44").Value = '0.7225'
$ws.Range("E44").Value = '  +1.38%  '
$ws.Range("D45").Value = '15.48'
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("D46").Value = '2.469'
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").Value = '4.176'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("D48").Value = '0.9972'
$ws.Range("E48").Value = '  -0.64%  '
$ws.Range("D49").Value = '141.95'
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").Value = '1.312'
$ws.Range("E50").Value = '  +13.93%  '
$ws.Range("D51").Value = '0.00000000376'
$ws.Range("E51").Value = '  -0.35%  '

# Restore default style on the touched D-column cells so no stray number-format
# style index is left behind on cells that did not have one originally.
$ws.Range("D2:D51").Style = "Normal"
